# "Error Calculations and Plots" - re-roll which cells in column E/F are
# treated as missing data, and drop two rows (RM 232, SC 92) from the sheet.
#
# All row numbers below refer to the ORIGINAL (pre-delete) layout, so the
# value edits are applied first, then the two obsolete rows are deleted last
# (which naturally shifts everything below them up by one / two rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F re-roll for rows 2-25 (rows above the deleted ones; untouched by the shift) ---
$ws.Range("F6").Value  = 16.43   # RM 21:  was missing -> now has a value
$ws.Range("F8").Value  = ""      # RM 38:  had a value  -> now missing
$ws.Range("F12").Value = 17.45   # RM 81:  was missing -> now has a value
$ws.Range("F14").Value = ""      # RM 90:  had a value  -> now missing
$ws.Range("F17").Value = 17.78   # RM 116: was missing -> now has a value
$ws.Range("F18").Value = 18.35   # RM 120: was missing -> now has a value
$ws.Range("F19").Value = ""      # RM 125: had a value  -> now missing
$ws.Range("F20").Value = ""      # RM 134: had a value  -> now missing
$ws.Range("F23").Value = 16.48   # RM 140: was missing -> now has a value

# --- Column E/F re-roll for the rows that survive below the two deletions ---
# (original row -> will become SC 101 / SC 105 / SC 119 / SC 120 / SC 193 after delete)
$ws.Range("E29").Value = -10     # SC 101: was missing -> now has a value
$ws.Range("F29").Value = ""      # SC 101: had a value  -> now missing

$ws.Range("E30").Value = ""      # SC 105: had a value  -> now missing

$ws.Range("E31").Value = ""      # SC 119: had a value  -> now missing

$ws.Range("E32").Value = -5.7    # SC 120: was missing -> now has a value

$ws.Range("E34").Value = ""      # SC 193: had a value  -> now missing

# --- Drop the two rows that no longer appear in the sheet ---
# Row 26 "RM 232" is removed outright; after that shift, the former row 28
# "SC 92" has become row 27, so deleting row 27 next removes it too.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

Write-Output "edit complete"
